$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "/RME" from the steel (S/LFM+CDL) description line in the B2 cell.
$cell = $ws.Range("B2")
$currentText = $cell.Value2
$newText = $currentText -replace [regex]::Escape("8% S/LFM+CDL/RME/H:1"), "8% S/LFM+CDL/H:1"
$cell.Value = $newText

# Apply wrap-text formatting to the cell and grow the row to fit the long text.
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.6
